$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to stay text (some values look numeric,
# e.g. "548.96" or "0.0000146") so Excel does not auto-coerce them to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '62.928.74'
$ws.Range('E2').Value = '  -2.11%  '
$ws.Range('D3').Value = '2.680.04'
$ws.Range('E3').Value = '  -2.40%  '
$ws.Range('D5').Value = '548.96'
$ws.Range('E5').Value = '  -4.31%  '
$ws.Range('D6').Value = '157.21'
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('E8').Value = '  -2.13%  '
$ws.Range('E9').Value = '  -4.02%  '
$ws.Range('E10').Value = '  -2.82%  '
$ws.Range('D11').Value = '0.367'
$ws.Range('E11').Value = '  -4.76%  '
$ws.Range('D12').Value = '5.08'
$ws.Range('E12').Value = '  -12.23%  '
$ws.Range('D13').Value = '3.156.90'
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('D14').Value = '25.98'
$ws.Range('E14').Value = '  -4.02%  '
$ws.Range('D15').Value = '62.830.70'
$ws.Range('E15').Value = '  -1.87%  '
$ws.Range('D16').Value = '0.0000146'
$ws.Range('E16').Value = '  -3.34%  '
$ws.Range('D17').Value = '2.684.50'
$ws.Range('E17').Value = '  -2.62%  '
$ws.Range('D18').Value = '11.86'
$ws.Range('E18').Value = '  -2.14%  '
$ws.Range('D19').Value = '4.57'
$ws.Range('E19').Value = '  -5.31%  '
$ws.Range('D20').Value = '342.59'
$ws.Range('E20').Value = '  -3.87%  '
$ws.Range('E21').Value = '  -4.77%  '
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = '0.503'
$ws.Range('E23').Value = '  -4.13%  '
$ws.Range('D24').Value = '63.41'
$ws.Range('E24').Value = '  -2.28%  '
$ws.Range('E25').Value = '  -1.54%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = '8.09'
$ws.Range('E27').Value = '  -5.17%  '
$ws.Range('D28').Value = '0.0₃0852'
$ws.Range('E28').Value = '  -6.88%  '
$ws.Range('E29').Value = '  -1.97%  '
$ws.Range('E30').Value = '  -2.55%  '
$ws.Range('D31').Value = '7.01'
$ws.Range('E31').Value = '  -4.33%  '
$ws.Range('D32').Value = '166.95'
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').Value = '4.78'
$ws.Range('E34').Value = '  -3.18%  '
$ws.Range('D35').Value = '19.49'
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('D36').Value = '1.42'
$ws.Range('E36').Value = '  -5.00%  '
$ws.Range('D37').Value = '1.76'
$ws.Range('E37').Value = '  -3.51%  '
$ws.Range('D38').Value = '337.72'
$ws.Range('E38').Value = '  -3.22%  '
$ws.Range('D39').Value = '6.15'
$ws.Range('E39').Value = '  -3.43%  '
$ws.Range('D40').Value = '0.928'
$ws.Range('E40').Value = '  -6.80%  '
$ws.Range('E41').Value = '  -2.08%  '
$ws.Range('D42').Value = '3.93'
$ws.Range('E42').Value = '  -5.55%  '
$ws.Range('D43').Value = '20.25'
$ws.Range('E43').Value = '  -5.97%  '
$ws.Range('D44').Value = '20.65'
$ws.Range('E44').Value = '  -7.86%  '
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('D46').Value = '0.0559'
$ws.Range('E46').Value = '  -5.33%  '
$ws.Range('D47').Value = '0.999'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D49').Value = '0.0970'
$ws.Range('E49').Value = '  -3.92%  '
$ws.Range('D50').Value = '128.52'
$ws.Range('E50').Value = '  -5.26%  '
$ws.Range('D51').Value = '2.086.32'
$ws.Range('E51').Value = '  -2.28%  '

# Restore the default (unstyled) cell style so the saved file matches the
# original formatting (these cells carry no explicit style).
$dataRange.Style = "Normal"
